# Append the new daily row (2025-08-23) to the portfolio-updates sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 8

# Column A holds dates stored as literal text in this sheet (matching the
# existing rows), so force text formatting before writing the value to
# avoid Excel auto-converting the "YYYY-MM-DD" string into a date serial,
# then clear the temporary number format so the cell keeps the sheet's
# default (unstyled) look, just like the other data rows.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "2025-08-23"
$ws.Range("A" + $newRow).ClearFormats()

$ws.Range("B" + $newRow).Value = 58.5099983215332
$ws.Range("C" + $newRow).Value = 680.2999877929688
$ws.Range("D" + $newRow).Value = 319.1000061035156
